$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3089087
$ws.Range("J76").Value = 4632267
$ws.Range("L76").Value = 4632267
$ws.Range("N76").Value = -4632897
$ws.Range("H79").Value = 3089087
$ws.Range("J79").Value = 4632267
$ws.Range("L79").Value = 4632267
$ws.Range("N79").Value = -4634451
$ws.Range("H98").Value = 785.55554
$ws.Range("I98").Value = 785.55554
$ws.Range("K98").Value = 785.55554
$ws.Range("M98").Value = 712.44446
$ws.Range("H122").Value = 785.55554
$ws.Range("I122").Value = 785.55554
$ws.Range("K122").Value = 2356.66662
$ws.Range("M122").Value = 93.33338000000003
$ws.Range("H138").Value = 4325.143
$ws.Range("I138").Value = 4008.3333
$ws.Range("J138").Value = 4411.5454
$ws.Range("K138").Value = 12024.9999
$ws.Range("L138").Value = 13234.6362
$ws.Range("M138").Value = -6884.999899999999
$ws.Range("N138").Value = -23514.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1011.1667
$ws.Range("I2").Value = 981.55554
$ws.Range("K2").Value = 981.55554
$ws.Range("M2").Value = -868.55554
$ws.Range("H32").Value = 9471.597
$ws.Range("I32").Value = 6956.7017
$ws.Range("J32").Value = 23806.5
$ws.Range("K32").Value = 6956.7017
$ws.Range("L32").Value = 23806.5
$ws.Range("M32").Value = -6669.7017
$ws.Range("N32").Value = -24380.5
$ws.Range("H45").Value = 2460.7715
$ws.Range("I45").Value = 2143.96
$ws.Range("K45").Value = 2143.96
$ws.Range("M45").Value = -1766.96
$ws.Range("H88").Value = 200993.4
$ws.Range("I88").Value = 1203
$ws.Range("K88").Value = 1203
$ws.Range("M88").Value = -797
$ws.Range("H91").Value = 200993.4
$ws.Range("I91").Value = 1203
$ws.Range("K91").Value = 1203
$ws.Range("M91").Value = 201
$ws.Range("H97").Value = 872.5
$ws.Range("I97").Value = 872.5
$ws.Range("K97").Value = 872.5
$ws.Range("M97").Value = -376.5
$ws.Range("H116").Value = 1011.1667
$ws.Range("I116").Value = 981.55554
$ws.Range("K116").Value = 981.55554
$ws.Range("M116").Value = 1312.44446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1011.1667
$ws.Range("I3").Value = 981.55554
$ws.Range("K3").Value = 981.55554
$ws.Range("M3").Value = -867.55554
$ws.Range("H94").Value = 716.96
$ws.Range("I94").Value = 590.8461
$ws.Range("K94").Value = 590.8461
$ws.Range("M94").Value = -139.8461

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3613.4443
$ws.Range("J31").Value = 5590.033
$ws.Range("L31").Value = 5590.033
$ws.Range("N31").Value = -6180.033
$ws.Range("H34").Value = 3613.4443
$ws.Range("J34").Value = 5590.033
$ws.Range("L34").Value = 5590.033
$ws.Range("N34").Value = -5994.033
$ws.Range("H62").Value = 3649.9167
$ws.Range("I62").Value = 3422.111
$ws.Range("K62").Value = 3422.111
$ws.Range("M62").Value = -2798.111
$ws.Range("H65").Value = 3649.9167
$ws.Range("I65").Value = 3422.111
$ws.Range("K65").Value = 17110.555
$ws.Range("M65").Value = -13990.555
$ws.Range("H141").Value = 22480.666
$ws.Range("J141").Value = 25976.8
$ws.Range("L141").Value = 25976.8
$ws.Range("N141").Value = -36336.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 885.5714
$ws.Range("J17").Value = 999.8
$ws.Range("L17").Value = 2999.4
$ws.Range("N17").Value = -3337.4
$ws.Range("H34").Value = 701.6111
$ws.Range("J34").Value = 765.26666
$ws.Range("L34").Value = 2295.79998
$ws.Range("N34").Value = -2463.79998
$ws.Range("H39").Value = 3607.2727
$ws.Range("J39").Value = 3607.2727
$ws.Range("L39").Value = 10821.8181
$ws.Range("N39").Value = -11409.8181
$ws.Range("H55").Value = 2500
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 4500
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 13500
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -13854
$ws.Range("H92").Value = 25000770
$ws.Range("J92").Value = 1499.5
$ws.Range("L92").Value = 4498.5
$ws.Range("N92").Value = -6994.5
$ws.Range("H131").Value = 726.47
$ws.Range("J131").Value = 760.68134
$ws.Range("L131").Value = 2282.04402
$ws.Range("N131").Value = -12362.04402
$ws.Range("H136").Value = 3402.5
$ws.Range("I136").Value = 1517
$ws.Range("J136").Value = 4973.75
$ws.Range("K136").Value = 4551
$ws.Range("L136").Value = 14921.25
$ws.Range("M136").Value = 549
$ws.Range("N136").Value = -25121.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25003540
$ws.Range("I80").Value = 38464680
$ws.Range("K80").Value = 38464680
$ws.Range("M80").Value = -38463682
$ws.Range("H83").Value = 25003540
$ws.Range("I83").Value = 38464680
$ws.Range("K83").Value = 192323400
$ws.Range("M83").Value = -192318408

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4266.6665
$ws.Range("I7").Value = 3580
$ws.Range("K7").Value = 3580
$ws.Range("M7").Value = -3468
$ws.Range("H46").Value = 950.0833
$ws.Range("I46").Value = 920.1
$ws.Range("K46").Value = 920.1
$ws.Range("M46").Value = -732.1
$ws.Range("H68").Value = 2832.5
$ws.Range("I68").Value = 2399
$ws.Range("K68").Value = 2399
$ws.Range("M68").Value = -1650
$ws.Range("H71").Value = 2832.5
$ws.Range("I71").Value = 2399
$ws.Range("K71").Value = 11995
$ws.Range("M71").Value = -8251
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 787725.0600000001
$ws.Range("I122").Value = 1092651.5
$ws.Range("K122").Value = 3277954.5
$ws.Range("M122").Value = -3275504.5
$ws.Range("H126").Value = 4266.6665
$ws.Range("I126").Value = 3580
$ws.Range("K126").Value = 10740
$ws.Range("M126").Value = -8270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 19805
$ws.Range("J94").Value = 19805
$ws.Range("L94").Value = 19805
$ws.Range("N94").Value = -21607
$ws.Range("H126").Value = 2402.3333
$ws.Range("I126").Value = 1977.7778
$ws.Range("K126").Value = 5933.3334
$ws.Range("M126").Value = -3463.3334
